# Setting tab TC-01,02,09 added (27/10/25)
# Adds a new worksheet "tc009" at the end of the workbook (after the
# current last sheet, "tc002") and populates it with a small
# fieldname/datatype table, matching the data already used on the
# "Demo" sheet for the "datatype" / "Text Box" values.

$wb = $excel.ActiveWorkbook

# Add the new worksheet right after the current last sheet so it lands
# at the end of the tab strip and becomes the active sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "tc009"

# Populate the new sheet's data.
$newSheet.Range("A1").Value = "fieldname"
$newSheet.Range("B1").Value = "datatype"
$newSheet.Range("A2").Value = "Category"
$newSheet.Range("B2").Value = "Text Box"

# Match the recorded selection on the new tab.
$newSheet.Range("C6").Select() | Out-Null
